$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 44990
$ws.Range("C12").Value = "2 hrs"
$ws.Range("D12").Value = "GITHUB issues and update with bug fix on text color "

$ws.Range("B17").Value = 44999
$ws.Range("C17").Value = "2 hrs"
$ws.Range("D17").Value = "adding splash activity and logos, icon for app but two icon shows not sure why have to check "

$ws.Range("B12").NumberFormat = "d-mmm"
$ws.Range("B17").NumberFormat = "d-mmm"
